$d = $word.ActiveDocument

# --- Change 1: append the red "(This is a change – Version for main branch)" text
# after the first paragraph's existing text, preceded by two spaces.
$findRange = $d.Content
$found = $findRange.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$findRange.Collapse(0)
$findRange.InsertAfter("  ")

$findRange.Collapse(0)
$findRange.InsertAfter("(This is a change – Ve")
$findRange.Font.Color = 255

$findRange.Collapse(0)
$findRange.InsertAfter("rsion for main branch")
$findRange.Font.Color = 255

$findRange.Collapse(0)
$findRange.InsertAfter(")")
$findRange.Font.Color = 255

# --- Change 2: add a new, empty paragraph with light-grey shading at the very
# end of the document (right before the section break).
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')
